$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

$ws.Range("B2").Value = 136
$ws.Range("B3").Value = 118
$ws.Range("B4").Value = 104
$ws.Range("B5").Value = 100
$ws.Range("B6").Value = 94
$ws.Range("B7").Value = 25
$ws.Range("B8").Value = 17
